# Add "Area" / "Atotal" columns (G, H) to the discharge worksheet,
# mirroring the existing "Q"/"Qtotal" (segment discharge) columns but for
# cross-sectional area instead of discharge.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Row 2: first segment area (uses 0 as the "previous" depth, matching the
# Q formula's special-cased first row) + running total via SUM(G2:G11)
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Rows 3-11: per-segment area, extending the table down to row 11
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"

$ws.Range("G1:H11").Select() | Out-Null

$excel.Calculate()
